$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Break Beam BOM")

$fiveMmName = "IR Break Beam Sensor with Premium Wire Header Ends - 5mm LEDs"
$fiveMmDatasheet = "IR Break Beam Sensor with Premium Wire Header Ends - 5mm LEDs : ID 2168 : Adafruit Industries, Unique & fun DIY electronics and kits"
$threeMmName = "IR Break Beam Sensors with Premium Wire Header Ends - 3mm LEDs"
$threeMmDatasheet = "IR Break Beam Sensors with Premium Wire Header Ends - 3mm LEDs : ID 2167 : Adafruit Industries, Unique & fun DIY electronics and kits"
$dimension = "(20*10*8)"

# Enter 3mm column (D) first to mimic original authoring order for shared-string indices
$ws.Range("D1").Value = $threeMmName
$ws.Range("D2").Value = $threeMmDatasheet
$ws.Range("D3").Value = $dimension

# Then 5mm column (B)
$ws.Range("B1").Value = $fiveMmName
$ws.Range("B2").Value = $fiveMmDatasheet
$ws.Range("B3").Value = $dimension

$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 1

$ws.Range("B5").Formula = "=10.24*5.95"
$ws.Range("D5").Formula = "=10.24*2.95"

$ws.Range("B6").Formula = "=B4*B5"

# Hyperlinks: D2 -> product 2167 (3mm), B2 -> product 2168 (5mm)
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.adafruit.com/product/2167", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.adafruit.com/product/2167") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.adafruit.com/product/2168", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.adafruit.com/product/2168") | Out-Null
